$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.512.80"
$ws.Range("E2").Value = "  +2.09%  "
$ws.Range("D3").Value = "1.681.68"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "217.67"
$ws.Range("E5").Value = "  +4.11%  "
$ws.Range("D6").Value = "0.5321"
$ws.Range("E6").Value = "  +2.92%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +4.44%  "
$ws.Range("D9").Value = "0.06424"
$ws.Range("E9").Value = "  +3.05%  "
$ws.Range("E10").Value = "  +5.58%  "
$ws.Range("D11").Value = "0.07806"
$ws.Range("E11").Value = "  +3.58%  "
$ws.Range("D12").Value = "1.687.23"
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("E13").Value = "  +3.44%  "
$ws.Range("D14").Value = "0.5619"
$ws.Range("E14").Value = "  +4.19%  "
$ws.Range("D15").Value = "0.0₅8418"
$ws.Range("E15").Value = "  +5.95%  "
$ws.Range("D16").Value = "66.03"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").Value = "26.550.18"
$ws.Range("E17").Value = "  +2.12%  "
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "4.801"
$ws.Range("E19").Value = "  +3.54%  "
$ws.Range("D20").Value = "196.01"
$ws.Range("E20").Value = "  +6.00%  "
$ws.Range("E21").Value = "  +3.90%  "
$ws.Range("D22").Value = "6.381"
$ws.Range("E22").Value = "  +4.80%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "143.17"
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").Value = "0.1280"
$ws.Range("E25").Value = "  +7.50%  "
$ws.Range("D26").Value = "7.467"
$ws.Range("E26").Value = "  +1.85%  "
$ws.Range("D27").Value = "16.19"
$ws.Range("E27").Value = "  +4.69%  "
$ws.Range("E28").Value = "  +3.65%  "
$ws.Range("D29").Value = "0.06132"
$ws.Range("E29").Value = "  +2.99%  "
$ws.Range("D30").Value = "1.279"
$ws.Range("E30").Value = "  +2.98%  "
$ws.Range("D31").Value = "3.608"
$ws.Range("E31").Value = "  +7.73%  "
$ws.Range("D32").Value = "3.458"
$ws.Range("E32").Value = "  +3.47%  "
$ws.Range("D33").Value = "1.707"
$ws.Range("E33").Value = "  +6.36%  "
$ws.Range("E34").Value = "  +4.83%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.419"
$ws.Range("E35").Value = "  +1.56%  "
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").Value = "2.789"
$ws.Range("E36").Value = "  +2.05%  "
$ws.Range("D37").Value = "0.5710"
$ws.Range("E37").Value = "  -2.23%  "
$ws.Range("D38").Value = "0.01642"
$ws.Range("E38").Value = "  +2.99%  "
$ws.Range("E39").Value = "  +4.04%  "
$ws.Range("D40").Value = "0.8715"
$ws.Range("E40").Value = "  +3.68%  "
$ws.Range("D41").Value = "1.060.58"
$ws.Range("E41").Value = "  +1.61%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "1.832.75"
$ws.Range("E44").Value = "  +2.19%  "
$ws.Range("D45").Value = "0.0₈112"
$ws.Range("E45").Value = "  +4.59%  "
$ws.Range("D46").Value = "57.28"
$ws.Range("E46").Value = "  +5.41%  "
$ws.Range("D47").Value = "8.144"
$ws.Range("E47").Value = "  +1.97%  "
$ws.Range("D48").Value = "0.9986"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "0.05202"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").Value = "6.078"
$ws.Range("E50").Value = "  +5.07%  "
$ws.Range("D51").Value = "0.4240"
$ws.Range("E51").Value = "  +0.18%  "
